$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.70915780501259
$ws.Range("C2").Value = 10.70920021105154
$ws.Range("E2").Value = 24.80484117407459
$ws.Range("F2").Value = 38.9703355692271
$ws.Range("G2").Value = 3.609357570187385
$ws.Range("J2").Value = 7.554241624709104
$ws.Range("O2").Value = 18.49811662722985
$ws.Range("B3").Value = 14.9225889089125
$ws.Range("C3").Value = 10.07614640023442
$ws.Range("E3").Value = 24.46313151764528
$ws.Range("F3").Value = 38.81397716675223
$ws.Range("G3").Value = 3.611612679233842
$ws.Range("J3").Value = 7.588878807803155
$ws.Range("O3").Value = 18.66998579175946
$ws.Range("B4").Value = 14.41768862667285
$ws.Range("C4").Value = 9.6651455683797
$ws.Range("E4").Value = 24.25525372428084
$ws.Range("F4").Value = 38.73097637385831
$ws.Range("G4").Value = 3.613067472433015
$ws.Range("J4").Value = 7.611325524439433
$ws.Range("O4").Value = 18.78289541656396
$ws.Range("B5").Value = 14.20662669745182
$ws.Range("C5").Value = 9.492111557157806
$ws.Range("E5").Value = 24.17112098167006
$ws.Range("F5").Value = 38.70044397735072
$ws.Range("G5").Value = 3.613678011051287
$ws.Range("J5").Value = 7.620769851048864
$ws.Range("O5").Value = 18.83075393398789
$ws.Range("B6").Value = 14.17126602470314
$ws.Range("C6").Value = 9.463046404370692
$ws.Range("E6").Value = 24.15718848530195
$ws.Range("F6").Value = 38.69557338697577
$ws.Range("G6").Value = 3.613780461222794
$ws.Range("J6").Value = 7.62235603519241
$ws.Range("O6").Value = 18.83881209291801
$ws.Range("B7").Value = 14.4148633754037
$ws.Range("C7").Value = 9.662834342476211
$ws.Range("E7").Value = 24.25411661347384
$ws.Range("F7").Value = 38.73055125376742
$ws.Range("G7").Value = 3.613075634631434
$ws.Range("J7").Value = 7.611451690224419
$ws.Range("O7").Value = 18.78353338705759
$ws.Range("B8").Value = 15.44261659681623
$ws.Range("C8").Value = 10.49558383003306
$ws.Range("E8").Value = 24.68667439726998
$ws.Range("F8").Value = 38.9137384997037
$ws.Range("G8").Value = 3.610120608679125
$ws.Range("J8").Value = 7.565940083545336
$ws.Range("O8").Value = 18.55583954406638
$ws.Range("B9").Value = 17.27652240782981
$ws.Range("C9").Value = 11.94970783343662
$ws.Range("E9").Value = 25.54615831245221
$ws.Range("F9").Value = 39.37497531057026
$ws.Range("G9").Value = 3.60487966834991
$ws.Range("J9").Value = 7.486024344233226
$ws.Range("O9").Value = 18.16832194522659
$ws.Range("B10").Value = 18.50532417189012
$ws.Range("C10").Value = 12.90721767830671
$ws.Range("E10").Value = 26.17900956848825
$ws.Range("F10").Value = 39.77412158467931
$ws.Range("G10").Value = 3.601362954775072
$ws.Range("J10").Value = 7.432964217028337
$ws.Range("O10").Value = 17.92017364645741
$ws.Range("B11").Value = 19.03739309466579
$ws.Range("C11").Value = 13.31854801249046
$ws.Range("E11").Value = 26.46607611257662
$ws.Range("F11").Value = 39.96828513863473
$ws.Range("G11").Value = 3.599834773047865
$ws.Range("J11").Value = 7.410046163813043
$ws.Range("O11").Value = 17.81536039372199
$ws.Range("B12").Value = 19.2349271327002
$ws.Range("C12").Value = 13.47081252208667
$ws.Range("E12").Value = 26.57457053468903
$ws.Range("F12").Value = 40.0435708930477
$ws.Range("G12").Value = 3.599266322556981
$ws.Range("J12").Value = 7.401542506545685
$ws.Range("O12").Value = 17.77684303479241
$ws.Range("B13").Value = 19.192561351965
$ws.Range("C13").Value = 13.43817530212761
$ws.Range("E13").Value = 26.5512151819687
$ws.Range("F13").Value = 40.02727931313575
$ws.Range("G13").Value = 3.599388294051657
$ws.Range("J13").Value = 7.403366148503783
$ws.Range("O13").Value = 17.78508605586189
$ws.Range("B14").Value = 19.05372387321936
$ws.Range("C14").Value = 13.33114504223073
$ws.Range("E14").Value = 26.4750067770453
$ws.Range("F14").Value = 39.974443966861
$ws.Range("G14").Value = 3.599787801401332
$ws.Range("J14").Value = 7.409343059648021
$ws.Range("O14").Value = 17.8121679620964
$ws.Range("B15").Value = 18.96816559205265
$ws.Range("C15").Value = 13.26513022099864
$ws.Range("E15").Value = 26.42829666575933
$ws.Range("F15").Value = 39.94230846831401
$ws.Range("G15").Value = 3.600033842931299
$ws.Range("J15").Value = 7.413026856903191
$ws.Range("O15").Value = 17.82890957127999
$ws.Range("B16").Value = 18.47000350998607
$ws.Range("C16").Value = 12.87984734292605
$ws.Range("E16").Value = 26.1602245319817
$ws.Range("F16").Value = 39.76168168324623
$ws.Range("G16").Value = 3.601464260833963
$ws.Range("J16").Value = 7.434486464476228
$ws.Range("O16").Value = 17.92718690755662
$ws.Range("B17").Value = 18.15744201297592
$ws.Range("C17").Value = 12.63726948414282
$ws.Range("E17").Value = 25.99549344822488
$ws.Range("F17").Value = 39.65406316759623
$ws.Range("G17").Value = 3.602360071842229
$ws.Range("J17").Value = 7.447963220534369
$ws.Range("O17").Value = 17.98955277699172
$ws.Range("B18").Value = 17.97513471465044
$ws.Range("C18").Value = 12.49546327974022
$ws.Range("E18").Value = 25.90067267484241
$ws.Range("F18").Value = 39.59335213756758
$ws.Range("G18").Value = 3.602882060232044
$ws.Range("J18").Value = 7.455829481121758
$ws.Range("O18").Value = 18.02618284215243
$ws.Range("B19").Value = 17.91297666831917
$ws.Range("C19").Value = 12.44705864363765
$ws.Range("E19").Value = 25.86855854737199
$ws.Range("F19").Value = 39.5730019793117
$ws.Range("G19").Value = 3.603059956107762
$ws.Range("J19").Value = 7.458512589510109
$ws.Range("O19").Value = 18.03871511258476
$ws.Range("B20").Value = 18.19097719743263
$ws.Range("C20").Value = 12.66332856742661
$ws.Range("E20").Value = 26.01303747092347
$ws.Range("F20").Value = 39.66539670902116
$ws.Range("G20").Value = 3.602264013904789
$ws.Range("J20").Value = 7.446516719882867
$ws.Range("O20").Value = 17.98283519400519
$ws.Range("B21").Value = 19.09461152412665
$ws.Range("C21").Value = 13.36267740332504
$ws.Range("E21").Value = 26.49739750241487
$ws.Range("F21").Value = 39.98991563624503
$ws.Range("G21").Value = 3.599670178973462
$ws.Range("J21").Value = 7.407582752146158
$ws.Range("O21").Value = 17.8041814083125
$ws.Range("B22").Value = 19.6621456410901
$ws.Range("C22").Value = 13.79935527087289
$ws.Range("E22").Value = 26.81268100380075
$ws.Range("F22").Value = 40.21224148807349
$ws.Range("G22").Value = 3.598034611652644
$ws.Range("J22").Value = 7.383156492182333
$ws.Range("O22").Value = 17.69426498581077
$ws.Range("B23").Value = 19.36137235995398
$ws.Range("C23").Value = 13.56815983565708
$ws.Range("E23").Value = 26.64455494911029
$ws.Range("F23").Value = 40.09266317313405
$ws.Range("G23").Value = 3.598902104694531
$ws.Range("J23").Value = 7.396100111157048
$ws.Range("O23").Value = 17.75229892978244
$ws.Range("B24").Value = 18.17582406882164
$ws.Range("C24").Value = 12.65155455914747
$ws.Range("E24").Value = 26.00510616831335
$ws.Range("F24").Value = 39.66026919681712
$ws.Range("G24").Value = 3.602307419966819
$ws.Range("J24").Value = 7.447170314303625
$ws.Range("O24").Value = 17.98586979972016
$ws.Range("B25").Value = 16.80074925015739
$ws.Range("C25").Value = 11.57570084204496
$ws.Range("E25").Value = 25.31300969901667
$ws.Range("F25").Value = 39.23945593487018
$ws.Range("G25").Value = 3.606238585563084
$ws.Range("J25").Value = 7.506648272792525
$ws.Range("O25").Value = 18.26678015199523
